$wb = $excel.ActiveWorkbook

# Remove Sheet2 and Sheet3 (workbook now only has Sheet1)
$wb.Worksheets.Item("Sheet2").Delete()
$wb.Worksheets.Item("Sheet3").Delete()

$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Append two new rows of measurements below the existing data
$ws.Range("C47").Value = 21
$ws.Range("D47").Value = 132
$ws.Range("E47").Value = 65

$ws.Range("C48").Value = 21.3
$ws.Range("D48").Value = 133
$ws.Range("E48").Value = 66

# Leave the selection on the newly scrolled-to area
$ws.Range("E43").Select()
